$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Lost Time" row values (row 4): B4 and C4
$ws.Range("B4").Value = 10
$ws.Range("C4").Value = 5

# Recalculate so the shared formula in D4 (B4-C4) updates
$excel.Calculate()

# Update the active selection to match the new cursor position
$ws.Range("C7").Select()
